$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the English-header row (rank/team/session/score/BOSS) - shifts rows up
$ws.Rows.Item(2).Delete()

# Remove the trailing "尤文图斯" row (now row 4 after the shift above) - shifts rows up
$ws.Rows.Item(4).Delete()

# Header row: "老板" -> "BOSS", and drop the bold/bordered header style entirely
$ws.Range("A1:E1").ClearFormats()
$ws.Range("E1").Value = "BOSS"

# F1 no longer used at all - remove it completely
$ws.Range("F1").Clear()

# New data cell introduced in row 2
$ws.Range("F2").Value = "test"

# Move the active selection off the old F1:F5 range
$ws.Range("H15").Select() | Out-Null
